$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Characters(21, 2).Text = "17"
$ws.Range("C9").Characters(27, 9).Text = "4/24/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/30/2023"

# --- Style-change cells: copy format from a stable donor cell, then set value ---
# Donor cells (unaffected by this edit):
#   D14 = style 14 + shared string "0"   (index 20)
#   E14 = style 14 + shared string "***.*" (index 21)
#   F15 = style 16 (integer count format)
#   L14 = style 15 (decimal/percent format)

$ws.Range("D14").Copy($ws.Range("C15"))
$ws.Range("D14").Copy($ws.Range("C16"))
$ws.Range("D14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("F15").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("L14").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -50
$ws.Range("D14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("D14").Copy($ws.Range("C23"))
$ws.Range("D14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("D14").Copy($ws.Range("C26"))
$ws.Range("D14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("L14").Copy($ws.Range("L28"))
$ws.Range("L28").Value = 100
$ws.Range("L14").Copy($ws.Range("L29"))
$ws.Range("L29").Value = 100

# --- Plain value updates (style unchanged) ---
$ws.Range("L15").Value = -50
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 24
$ws.Range("K16").Value = 4.347826086956
$ws.Range("L16").Value = 26.315789473684
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 78
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = -2.5
$ws.Range("L17").Value = 69.565217391304
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 31
$ws.Range("K18").Value = -48.387096774193
$ws.Range("L18").Value = -46.666666666666
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 104
$ws.Range("J19").Value = 159
$ws.Range("K19").Value = -34.591194968553
$ws.Range("L19").Value = 13.043478260869
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 9
$ws.Range("H20").Value = 125
$ws.Range("I20").Value = 29
$ws.Range("K20").Value = 16
$ws.Range("L20").Value = 141.666666666667
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 6.25
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 255
$ws.Range("J21").Value = 321
$ws.Range("K21").Value = -20.560747663551
$ws.Range("L21").Value = 22.009569377990
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 178.571428571429
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 12.121212121212
$ws.Range("I24").Value = 405
$ws.Range("J24").Value = 398
$ws.Range("K24").Value = 1.758793969849
$ws.Range("L24").Value = 30.225080385852
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 13
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -14.545454545454
$ws.Range("I25").Value = 183
$ws.Range("J25").Value = 162
$ws.Range("K25").Value = 12.962962962963
$ws.Range("L25").Value = 72.641509433962
$ws.Range("L26").Value = -33.333333333333
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = 76.923076923076
$ws.Range("L27").Value = 27.777777777777
